# Generate Report for Handback
# Updates handoff/handback timestamps and priority (ht -> mt) to reflect a
# newer report generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
# 2016-08-17 08:14:45 -> 2016-08-17 08:15:51
$wsOverview.Range("G2").Value = "2016-08-17 08:15:51"
$wsOverview.Range("G3").Value = "2016-08-17 08:15:51"

# --- zh-cn sheet ---
# Priority column (E): ht -> mt
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"

# Correspond Handoff Datetime (H): 2016-08-17 08:14:39 -> 2016-08-17 08:15:46
$wsZhCn.Range("H2").Value = "2016-08-17 08:15:46"
$wsZhCn.Range("H3").Value = "2016-08-17 08:15:46"

# Correspond Handback DateTime (K): 2016-08-17 08:15:19 -> 2016-08-17 08:16:19
$wsZhCn.Range("K2").Value = "2016-08-17 08:16:19"
$wsZhCn.Range("K3").Value = "2016-08-17 08:16:19"

# --- de-de sheet ---
# Priority column (E): ht -> mt
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"

# Correspond Handoff Datetime (H): 2016-08-17 08:14:45 -> 2016-08-17 08:15:51
$wsDeDe.Range("H2").Value = "2016-08-17 08:15:51"
$wsDeDe.Range("H3").Value = "2016-08-17 08:15:51"

# Correspond Handback DateTime (K): 2016-08-17 08:15:26 -> 2016-08-17 08:16:27
$wsDeDe.Range("K2").Value = "2016-08-17 08:16:27"
$wsDeDe.Range("K3").Value = "2016-08-17 08:16:27"

$wb.Save()
